# Build Room / Log update script
# - Updates LastCount (B) / NewCount (C) figures on the "Build Room" sheet (Sheet1)
#   to reflect the latest stock movements.
# - Appends the corresponding movement rows to the "Log" sheet (Sheet2).
# - Normalizes the style of the previously-last log row (row 10) and the newly
#   added rows (11-27) to the sheet's standard "Normal" cell style, leaving the
#   brand-new final row (28) with the default style (matching how the log has
#   historically been appended).

$wb = $excel.ActiveWorkbook
$wsBuild = $wb.Worksheets.Item("Sheet1")
$wsLog   = $wb.Worksheets.Item("Sheet2")

# Special characters used in a couple of item names (curly right double quote
# and a non-breaking space), built explicitly so they match the characters
# already used for these items elsewhere in the workbook.
$curlyQuote = [char]0x201D
$nbsp       = [char]0x00A0
$monitor24Name = "Monitor 24" + $curlyQuote + $nbsp
$wirelessHeadsetName = "Wireless Headset" + $nbsp + "Poly "

# ---------------------------------------------------------------------------
# 1. Build Room sheet - refresh LastCount / NewCount values per item
# ---------------------------------------------------------------------------
$buildUpdates = @(
    @{ Row = 2;  B = 147; C = 152 },   # Desktop Mini
    @{ Row = 3;  B = -1;  C = 399 },   # Dock Thunderbolt G4
    @{ Row = 4;  B = 60;  C = 40  },   # Laptop 840 G10
    @{ Row = 5;  B = 12;  C = 17  },   # Laptop 840 G9
    @{ Row = 7;  B = $null; C = 51 }, # Laptop x360 G8 (LastCount unchanged)
    @{ Row = 8;  B = 1;   C = 11  },   # Monitor 24"
    @{ Row = 10; B = 46;  C = 56  },   # USB External DVD-RW Drive
    @{ Row = 11; B = 30;  C = 40  },   # Wired Headset Poly 3325
    @{ Row = 13; B = 6;   C = 16  },   # Wireless Headset Poly
    @{ Row = 14; B = 5;   C = 15  }    # Wireless Keyboard and Mouse
)

foreach ($u in $buildUpdates) {
    if ($null -ne $u.B) {
        $wsBuild.Cells.Item($u.Row, 2).Value = $u.B
    }
    $wsBuild.Cells.Item($u.Row, 3).Value = $u.C
}

# ---------------------------------------------------------------------------
# 2. Log sheet - normalize style of row 10 (previously the last log entry)
# ---------------------------------------------------------------------------
$wsLog.Range("A10:C10").Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Log sheet - append the new movement log entries (rows 11-28)
# ---------------------------------------------------------------------------
$logRows = @(
    @{ Row = 11; Timestamp = "2023-12-05 21:12:00"; Item = "Desktop Mini"; Action = "Subtract 1" },
    @{ Row = 12; Timestamp = "2023-12-05 21:13:29"; Item = "Laptop x360 G8"; Action = "Subtract 2" },
    @{ Row = 13; Timestamp = "2023-12-05 21:13:52"; Item = "Dock Thunderbolt G4"; Action = "Add 1" },
    @{ Row = 14; Timestamp = "2023-12-05 21:58:55"; Item = "Laptop 840 G9"; Action = "Add 10" },
    @{ Row = 15; Timestamp = "2023-12-05 21:59:03"; Item = "Desktop Mini"; Action = "Add 10" },
    @{ Row = 16; Timestamp = "2023-12-05 21:59:09"; Item = $monitor24Name; Action = "Add 10" },
    @{ Row = 17; Timestamp = "2023-12-05 21:59:12"; Item = "Wired Headset Poly 3325"; Action = "Add 10" },
    @{ Row = 18; Timestamp = "2023-12-05 21:59:15"; Item = "USB External DVD-RW Drive"; Action = "Add 10" },
    @{ Row = 19; Timestamp = "2023-12-05 21:59:19"; Item = "Wireless Keyboard and Mouse"; Action = "Add 10" },
    @{ Row = 20; Timestamp = "2023-12-05 21:59:22"; Item = $wirelessHeadsetName; Action = "Add 10" },
    @{ Row = 21; Timestamp = "2023-12-07 22:09:25"; Item = "Dock Thunderbolt G4"; Action = "Subtract 2" },
    @{ Row = 22; Timestamp = "2023-12-07 22:09:39"; Item = "Dock Thunderbolt G4"; Action = "Subtract 20" },
    @{ Row = 23; Timestamp = "2023-12-07 22:09:53"; Item = "Laptop 840 G10"; Action = "Subtract 20" },
    @{ Row = 24; Timestamp = "2023-12-07 22:09:59"; Item = "Laptop 840 G10"; Action = "Subtract 20" },
    @{ Row = 25; Timestamp = "2023-12-07 22:10:03"; Item = "Laptop 840 G10"; Action = "Subtract 20" },
    @{ Row = 26; Timestamp = "2023-12-07 22:20:13"; Item = "Desktop Mini"; Action = "Add 5" },
    @{ Row = 27; Timestamp = "2023-12-07 22:20:17"; Item = "Laptop 840 G9"; Action = "Add 5" },
    @{ Row = 28; Timestamp = "2023-12-07 22:44:46"; Item = "Dock Thunderbolt G4"; Action = "Add 400" }
)

foreach ($entry in $logRows) {
    $wsLog.Cells.Item($entry.Row, 1).Value = $entry.Timestamp
    $wsLog.Cells.Item($entry.Row, 2).Value = $entry.Item
    $wsLog.Cells.Item($entry.Row, 3).Value = $entry.Action

    # Rows 11-27 adopt the sheet's "Normal" style explicitly (like row 10 above);
    # the brand-new last row (28) is left with the default/no explicit style.
    if ($entry.Row -le 27) {
        $wsLog.Range("A" + $entry.Row + ":C" + $entry.Row).Style = "Normal"
    }
}
